# Weekly price update: a new "Zapallo italiano" record (week of 2023-01-06)
# is inserted at row 401, pushing every existing record from row 401 down
# by one (to row 498). The sheet's dimension grows from A1:R497 to A1:R498.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 401; everything below shifts
# down one row (matches native Excel "Insert Sheet Rows" behaviour,
# including carrying the date-formatted style down from the row above).
$ws.Rows("401:401").Insert()

# Populate the newly inserted row 401 with the new weekly record.
$ws.Range("A401").Value = 9
$ws.Range("B401").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C401").Value = "Metropolitana"
$ws.Range("D401").Value = 44932
$ws.Range("E401").Value = 13
$ws.Range("F401").Value = 100112032
$ws.Range("G401").Value = "Zapallo italiano"
$ws.Range("H401").Value = "Sin especificar"
$ws.Range("I401").Value = "Primera"
$ws.Range("J401").Value = 340
$ws.Range("K401").Value = 4000
$ws.Range("L401").Value = 5000
$ws.Range("M401").Value = 4500
$ws.Range("N401").Value = "`$/caja 50 unidades"
$ws.Range("O401").Value = "Región Metropolitana"
$ws.Range("P401").Value = 90
$ws.Range("Q401").Value = 50
$ws.Range("R401").Value = "Hortaliza"
